$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value (Column C "Course General Info Sentence Count" updates,
# plus B60 which also changed)
$updates = @{
    "C3"  = 10
    "C5"  = 15
    "C8"  = 13
    "C10" = 10
    "C13" = 13
    "C15" = 12
    "C16" = 10
    "C17" = 13
    "C18" = 12
    "C19" = 10
    "C21" = 15
    "C22" = 10
    "C23" = 20
    "C27" = 18
    "C28" = 15
    "C31" = 10
    "C33" = 16
    "C40" = 14
    "C45" = 12
    "C47" = 14
    "C49" = 15
    "C51" = 16
    "C52" = 14
    "C53" = 16
    "C55" = 30
    "C56" = 16
    "C59" = 22
    "B60" = 20
    "C60" = 23
    "C61" = 17
    "C64" = 15
    "C68" = 13
    "C69" = 12
    "C70" = 13
    "C72" = 10
    "C78" = 14
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
